$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "Đánh giá (Incentive 6 tháng ...)" columns (G:H) from the
# payroll "Data" sheet. This shifts all subsequent columns left by two and
# adjusts the sheet dimension/row spans/shared formulas accordingly.
$ws.Columns("G:H").Delete()

# Leave the selection on the former last column (now O1), matching the
# post-edit view state.
$ws.Range("O1").Select()
